$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trends Status")

# Rapid Decline row
$ws.Range("C2").Value = 2
$ws.Range("E2").Value = 50

# Decline row
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0

# Stable row
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 50
